# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" column (E16:E30) listed the account-statement periods
# in ascending order (1811 .. 2001). The database refresh flips that list so
# the newest period (2001) is now on top and the oldest (1811) on bottom -
# i.e. old statements are "removed" from the top and new ones "added",
# pushing the list in descending order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @("2001","1912","1911","1910","1909","1908","1907","1906","1905","1904","1903","1902","1901","1812","1811")

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value2 = $periodos[$i]
}
